# Rename the inline logo pictures that live in the document's headers and
# footers:
#   - The BTEC "Logo-Orange" picture (in both headers)  image2.jpg -> image1.jpg
#   - The Pearson logo picture       (in both footers)  image1.png -> image2.png
#
# The pictures are matched by their (stable) alternative text / description
# so the correct shape is renamed regardless of header/footer ordering.
# The shape is selected first and renamed through the Selection object,
# which is the reliable path for InlineShapes that live in header/footer
# stories.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    foreach ($hf in $sec.Headers) {
        if ($hf.Exists) {
            $count = $hf.Range.InlineShapes.Count
            for ($i = 1; $i -le $count; $i++) {
                $shp = $hf.Range.InlineShapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Range.Select()
                    $word.Selection.InlineShapes.Item(1).Name = "image1.jpg"
                }
            }
        }
    }

    foreach ($hf in $sec.Footers) {
        if ($hf.Exists) {
            $count = $hf.Range.InlineShapes.Count
            for ($i = 1; $i -le $count; $i++) {
                $shp = $hf.Range.InlineShapes.Item($i)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Range.Select()
                    $word.Selection.InlineShapes.Item(1).Name = "image2.png"
                }
            }
        }
    }
}
